$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2, 3, 4) get rotated: new row2 <- old row3, new row3 <- old row4,
# new row4 <- old row2 (for columns A, B, E, F, G, H, Q, R, Z, AB).
# Additionally, the Q and R (Ost/Nord coordinate) values get rounded to whole numbers
# in the process.

$cols = @("A","B","E","F","G","H","Q","R","Z","AB")

# Capture current values for columns A, B, E, F, G, H, Q, R, Z, AB for rows 2-4
$orig = @{}
foreach ($r in 2..4) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowVals
}

# New row 2 gets old row 3 data, new row 3 gets old row 4 data, new row 4 gets old row 2 data
$mapping = @{ 2 = 3; 3 = 4; 4 = 2 }

foreach ($destRow in 2..4) {
    $srcRow = $mapping[$destRow]
    $srcVals = $orig[$srcRow]
    foreach ($c in $cols) {
        $val = $srcVals[$c]
        if ($c -eq "Q" -or $c -eq "R") {
            $val = [Math]::Round([double]$val)
        }
        $ws.Range("$c$destRow").Value2 = $val
    }
}
